# Remove the "Why Go!" slide from the deck.
#
# The deck originally had 18 slides; slide #3 ("Why Go!") is removed,
# leaving the remaining 17 slides in their original relative order
# (slide ids/relationship ids shift down to fill the gap, which
# PowerPoint / the COM host handles automatically on Delete()).

$p = $ppt.ActivePresentation

$target = $null
for ($i = 1; $i -le $p.Slides.Count; $i++) {
    $slide = $p.Slides.Item($i)
    $titleText = ""
    foreach ($shape in $slide.Shapes) {
        if ($shape.HasTextFrame -and $shape.TextFrame.HasText) {
            $titleText = $shape.TextFrame.TextRange.Text
            break
        }
    }
    if ($titleText -eq "Why Go!") {
        $target = $slide
        break
    }
}

if ($target -ne $null) {
    $target.Delete()
} else {
    # Fallback: the slide is known to be the 3rd slide in the original deck.
    $p.Slides.Item(3).Delete()
}
